# Generate Report for Handoff
# Updates the localization status report to reflect a fresh handoff:
#  - Status text "Handed back: in sync with en-US" -> "Ready for handoff"
#  - Refreshed timestamps for the handoff generation / handoff datetime
#  - Narrower width for the "Status" related columns on all sheets

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Status text updates ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# --- Timestamp updates ---
$overview.Range("G2").Value = "2016-08-31 15:17:38"
$dede.Range("H2").Value = "2016-08-31 15:17:38"
$zhcn.Range("H2").Value = "2016-08-31 15:17:33"

# --- Column width updates ---
# (target character width ~17.216 - closest the engine's pixel-quantized
# ColumnWidth setter can reach is 17.1666..; 16.36 maps to it reliably)
$overview.Range("E1").ColumnWidth = 16.36
$overview.Range("F1").ColumnWidth = 16.36
$zhcn.Range("C1").ColumnWidth = 16.36
$dede.Range("C1").ColumnWidth = 16.36
